$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1486
$ws.Range("G2").Value = 50
$ws.Range("F4").Value = 1778
$ws.Range("F7").Value = 662
$ws.Range("F9").Value = 65
$ws.Range("F12").Value = 4
$ws.Range("F13").Value = 82
$ws.Range("F14").Value = 153
$ws.Range("F15").Value = 23
$ws.Range("F16").Value = 123
$ws.Range("F17").Value = 74
$ws.Range("F18").Value = 106
$ws.Range("F19").Value = 4923
$ws.Range("F21").Value = 829
$ws.Range("F22").Value = 113
$ws.Range("F23").Value = 2233
$ws.Range("F25").Value = 20
$ws.Range("F26").Value = 2084

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 79

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1486
$ws.Range("G2").Value = 50
$ws.Range("F4").Value = 1778
$ws.Range("F7").Value = 662
$ws.Range("F9").Value = 65
$ws.Range("F12").Value = 4
$ws.Range("F13").Value = 82
$ws.Range("F14").Value = 153
$ws.Range("F15").Value = 23
$ws.Range("F16").Value = 123
$ws.Range("F17").Value = 74
$ws.Range("F18").Value = 106
$ws.Range("F19").Value = 4923
$ws.Range("F20").Value = 79
$ws.Range("F23").Value = 829
$ws.Range("F24").Value = 113
$ws.Range("F25").Value = 2233
$ws.Range("F27").Value = 20
$ws.Range("F28").Value = 2084
